$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - only B2 changes
$ws.Range("B2").Value = 10223522377862740

# Row 3 (RandomForestRegressor) - B3, C3, D3 change
$ws.Range("B3").Value = 0.01089868156067072
$ws.Range("C3").Value = 0.01232332837280409
$ws.Range("D3").Value = 5143037499490696

# Row 4: label GradientBoostingRegressor -> DecisionTreeRegressor, B4, C4, D4 change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01244444200575735
$ws.Range("C4").Value = 0.01400735442123902
$ws.Range("D4").Value = 0.02330753297411731

# Row 5: label AdaBoostRegressor -> MLPRegressor, B5, C5, D5 change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 373981365044435.4
$ws.Range("C5").Value = 203847449768752.8
$ws.Range("D5").Value = 793492095635976.4
